$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with two new columns P and Q, reusing the existing
# header style (bold, centered, thin border) from column O1.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For each data row (2..25), swap I/K and M/O values, and append P,Q = 2,2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P -> 2
    $ws.Cells.Item($r, 17).Value = 2   # Q -> 2
}
